$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '''29.971.32'
$ws.Range("E2").Value = '  -0.28%  '
$ws.Range("D3").Value = '''1.896.53'
$ws.Range("E3").Value = '  -0.87%  '
$ws.Range("D4").Value = '''1.002'
$ws.Range("E4").Value = '  +0.15%  '
$ws.Range("D5").Value = '''0.8359'
$ws.Range("E5").Value = '  +2.94%  '
$ws.Range("D6").Value = '''241.86'
$ws.Range("E6").Value = '  -0.47%  '
$ws.Range("D7").Value = '''1.001'
$ws.Range("E7").Value = '  +0.09%  '
$ws.Range("D8").Value = '''0.3288'
$ws.Range("E8").Value = '  +1.70%  '
$ws.Range("D9").Value = '''26.55'
$ws.Range("E9").Value = '  +0.49%  '
$ws.Range("D10").Value = '''0.07045'
$ws.Range("E10").Value = '  +1.15%  '
$ws.Range("D11").Value = '''0.08067'
$ws.Range("E11").Value = '  +0.46%  '
$ws.Range("D12").Value = '''0.7593'
$ws.Range("E12").Value = '  +0.99%  '
$ws.Range("D13").Value = '''1.899.54'
$ws.Range("E13").Value = '  -0.30%  '
$ws.Range("D14").Value = '''5.246'
$ws.Range("E14").Value = '  +0.08%  '
$ws.Range("D15").Value = '''92.17'
$ws.Range("E15").Value = '  -1.58%  '
$ws.Range("D16").Value = '''29.985.18'
$ws.Range("E16").Value = '  -0.26%  '
$ws.Range("E17").Value = '  +0.30%  '
$ws.Range("D18").Value = '''5.863'
$ws.Range("E18").Value = '  -2.39%  '
$ws.Range("E19").Value = '  -1.94%  '
$ws.Range("D20").Value = '''0.000007762'
$ws.Range("E20").Value = '  -0.37%  '
$ws.Range("E21").Value = '  +0.11%  '
$ws.Range("D22").Value = '''2.149.69'
$ws.Range("E22").Value = '  -0.26%  '
$ws.Range("D23").Value = '''1.002'
$ws.Range("E23").Value = '  +0.16%  '
$ws.Range("D24").Value = '''6.968'
$ws.Range("E24").Value = '  -0.41%  '
$ws.Range("D25").Value = '''0.1736'
$ws.Range("E25").Value = '  +21.40%  '
$ws.Range("B26").Value = 'Cosmos'
$ws.Range("C26").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range("D26").Value = '''9.242'
$ws.Range("E26").Value = '  -0.72%  '
$ws.Range("B27").Value = 'Monero'
$ws.Range("C27").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D27").Value = '''165.94'
$ws.Range("E27").Value = '  -1.60%  '
$ws.Range("E28").Value = '  -0.82%  '
$ws.Range("E29").Value = '  +1.74%  '
$ws.Range("E30").Value = '  -2.19%  '
$ws.Range("D31").Value = '''1.514'
$ws.Range("E31").Value = '  -0.60%  '
$ws.Range("D32").Value = '''0.05896'
$ws.Range("E32").Value = '  +10.48%  '
$ws.Range("D33").Value = '''4.283'
$ws.Range("E33").Value = '  -1.88%  '
$ws.Range("D34").Value = '''4.072'
$ws.Range("E34").Value = '  -1.43%  '
$ws.Range("E35").Value = '  -0.11%  '
$ws.Range("D36").Value = '''0.7284'
$ws.Range("E36").Value = '  -1.54%  '
$ws.Range("D37").Value = '''2.723'
$ws.Range("E37").Value = '  -0.34%  '
$ws.Range("D38").Value = '''0.01921'
$ws.Range("E38").Value = '  -0.52%  '
$ws.Range("D39").Value = '''2.773'
$ws.Range("E39").Value = '  -0.52%  '
$ws.Range("D40").Value = '''0.4426'
$ws.Range("E40").Value = '  -1.17%  '
$ws.Range("D41").Value = '''72.43'
$ws.Range("E41").Value = '  -0.47%  '
$ws.Range("B42").Value = 'TrustWalletToken'
$ws.Range("C42").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range("D42").Value = '''0.8577'
$ws.Range("E42").Value = '  +3.19%  '
$ws.Range("B43").Value = 'FraxShare'
$ws.Range("C43").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D43").Value = '''5.858'
$ws.Range("E43").Value = '  -4.99%  '
$ws.Range("E44").Value = '  +0.00%  '
$ws.Range("E45").Value = '  -1.13%  '
$ws.Range("D46").Value = '''101.80'
$ws.Range("E46").Value = '  +1.00%  '
$ws.Range("D47").Value = '''1.009.82'
$ws.Range("E47").Value = '  +4.18%  '
$ws.Range("B48").Value = 'EnergySwap'
$ws.Range("C48").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D48").Value = '''9.836'
$ws.Range("E48").Value = '  +0.13%  '
$ws.Range("B49").Value = 'Aptos'
$ws.Range("C49").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D49").Value = '''7.556'
$ws.Range("E49").Value = '  -0.96%  '
$ws.Range("D50").Value = '''2.045.12'
$ws.Range("D51").Value = '''35.85'
$ws.Range("E51").Value = '  -1.85%  '
